$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Row 1: 99.94 -> 0M
$t.Cell(1, 1).Range.Text = "0M"

# Row 2: 0.14 -> 0M
$t.Cell(2, 1).Range.Text = "0M"

# Row 3: 241 -> 0M
$t.Cell(3, 1).Range.Text = "0M"

# Row 4: 202 -> 807
$t.Cell(4, 1).Range.Text = "807"

# Row 6: 0.00022 -> 0.00058
$t.Cell(6, 1).Range.Text = "0.00058"

# Row 7: 0.00009 -> 0.00018
$t.Cell(7, 1).Range.Text = "0.00018"

# Row 8: 0.00003 -> 0.00006
$t.Cell(8, 1).Range.Text = "0.00006"

# Row 9: 0.00007 -> 0.00024
$t.Cell(9, 1).Range.Text = "0.00024"

# Row 10: 0.00008 -> 0.00026
$t.Cell(10, 1).Range.Text = "0.00026"

# Row 11: 0.00011 -> 0.00036
$t.Cell(11, 1).Range.Text = "0.00036"

# Row 12: 0.01754 -> 0.14423
$t.Cell(12, 1).Range.Text = "0.14423"

# Row 44: "202\t0.00015\t0.00058\t0.00030\t0.00009\t0.00024\t0.00026\t0.00036\t0.06155\t100.0" -> "99.94"
$t.Cell(44, 1).Range.Text = "99.94"

# Row 45: "201\t0.00009\t0.00037\t0.00016\t0.00005\t0.00013\t0.00016\t0.00018\t0.03242\t100.0" -> "0.14"
$t.Cell(45, 1).Range.Text = "0.14"

# Row 46: "202\t0.00008\t0.00034\t0.00016\t0.00005\t0.00012\t0.00013\t0.00021\t0.03272\t100.0" -> "241"
$t.Cell(46, 1).Range.Text = "241"
